$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 636733.2941520028
$ws.Range("C3").Value = 622606.6389124633
$ws.Range("C4").Value = 445948.2001064043
$ws.Range("C5").Value = 632173.8398814903
$ws.Range("C6").Value = 645430.2588331797
$ws.Range("C7").Value = 576964.7709690979
$ws.Range("C8").Value = 600611.1212153231
$ws.Range("C9").Value = 535593.3298955993
$ws.Range("C10").Value = 541483.8851404014
$ws.Range("C11").Value = 656266.418216069
$ws.Range("C12").Value = 613106.0944495659
$ws.Range("C13").Value = 578096.3282618701
$ws.Range("C14").Value = 712858.6685511762
